$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Help.xml")

# Remove the "Decision Trees" / "Clas-SLIQ" category block (rows 1-4, incl. leading blank row)
$ws2.Range("A1:A4").EntireRow.Delete()

# Remove "Clas-DMEL" and "Clas-GIL" rows from the "Genetic Rule Learning" block
$ws2.Range("A4:A5").EntireRow.Delete()

# Remove "Clas-DataSqueezer" and "Clas-Swap1" rows from the "Rule Learning" block
$ws2.Range("A12:A13").EntireRow.Delete()
$ws2.Range("A10:A11").EntireRow.Delete()

# Activate "Help.xml" sheet and set its selection, making it the active tab
$ws2.Activate()
$wb.Application.Goto($ws2.Range("A4:XFD5"))
